$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 9: set H=166.71428, I=166.71428, J=0, K=166.71428, L=0; clear M
$ws.Range("H9").Value = 166.71428
$ws.Range("I9").Value = 166.71428
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 166.71428
$ws.Range("L9").Value = 0
$ws.Range("M9").ClearContents()

# Row 32: set H=2122.6924, J=1863.2727, L=1863.2727, N=-2515.2727
$ws.Range("H32").Value = 2122.6924
$ws.Range("J32").Value = 1863.2727
$ws.Range("L32").Value = 1863.2727
$ws.Range("N32").Value = -2515.2727

# Row 74: set H=7088.7334, I=6043.857, K=6043.857, M=-5107.857
$ws.Range("H74").Value = 7088.7334
$ws.Range("I74").Value = 6043.857
$ws.Range("K74").Value = 6043.857
$ws.Range("M74").Value = -5107.857

# Row 77: set H=7088.7334, I=6043.857, K=30219.285, M=-25539.285
$ws.Range("H77").Value = 7088.7334
$ws.Range("I77").Value = 6043.857
$ws.Range("K77").Value = 30219.285
$ws.Range("M77").Value = -25539.285

# Row 86: set H=4619.3335, I=3082.8333, J=5387.5835, K=3082.8333, L=5387.5835, M=-1959.8333, N=-7633.5835
$ws.Range("H86").Value = 4619.3335
$ws.Range("I86").Value = 3082.8333
$ws.Range("J86").Value = 5387.5835
$ws.Range("K86").Value = 3082.8333
$ws.Range("L86").Value = 5387.5835
$ws.Range("M86").Value = -1959.8333
$ws.Range("N86").Value = -7633.5835

# Row 89: set H=4619.3335, I=3082.8333, J=5387.5835, K=15414.1665, L=26937.9175, M=-9798.166499999999, N=-38169.9175
$ws.Range("H89").Value = 4619.3335
$ws.Range("I89").Value = 3082.8333
$ws.Range("J89").Value = 5387.5835
$ws.Range("K89").Value = 15414.1665
$ws.Range("L89").Value = 26937.9175
$ws.Range("M89").Value = -9798.166499999999
$ws.Range("N89").Value = -38169.9175

# Row 113: set H=4040.6365, J=4999.4, L=4999.4, N=-11507.4
$ws.Range("H113").Value = 4040.6365
$ws.Range("J113").Value = 4999.4
$ws.Range("L113").Value = 4999.4
$ws.Range("N113").Value = -11507.4

# Row 116: set H=4298.4546, I=3766.4285, J=5229.5, K=3766.4285, L=5229.5, M=-324.4285, N=-12113.5
$ws.Range("H116").Value = 4298.4546
$ws.Range("I116").Value = 3766.4285
$ws.Range("J116").Value = 5229.5
$ws.Range("K116").Value = 3766.4285
$ws.Range("L116").Value = 5229.5
$ws.Range("M116").Value = -324.4285
$ws.Range("N116").Value = -12113.5

# Row 132: set H=1529.05, I=1114.3438, J=3187.875, K=3343.0314, L=9563.625, M=-813.0314000000003, N=-14623.625
$ws.Range("H132").Value = 1529.05
$ws.Range("I132").Value = 1114.3438
$ws.Range("J132").Value = 3187.875
$ws.Range("K132").Value = 3343.0314
$ws.Range("L132").Value = 9563.625
$ws.Range("M132").Value = -813.0314000000003
$ws.Range("N132").Value = -14623.625

# Row 137: set H=7125.3335, I=3146.1333, J=13757.333, K=9438.3999, L=41271.999, M=-6888.3999, N=-46371.999
$ws.Range("H137").Value = 7125.3335
$ws.Range("I137").Value = 3146.1333
$ws.Range("J137").Value = 13757.333
$ws.Range("K137").Value = 9438.3999
$ws.Range("L137").Value = 41271.999
$ws.Range("M137").Value = -6888.3999
$ws.Range("N137").Value = -46371.999

# Row 138: set H=3528.6667, I=1178.3125, J=5163.696, K=3534.9375, L=15491.088, M=1605.0625, N=-25771.088
$ws.Range("H138").Value = 3528.6667
$ws.Range("I138").Value = 1178.3125
$ws.Range("J138").Value = 5163.696
$ws.Range("K138").Value = 3534.9375
$ws.Range("L138").Value = 15491.088
$ws.Range("M138").Value = 1605.0625
$ws.Range("N138").Value = -25771.088

# Row 141: set H=1408.9584, I=1378.1818, K=4134.5454, M=1045.4546
$ws.Range("H141").Value = 1408.9584
$ws.Range("I141").Value = 1378.1818
$ws.Range("K141").Value = 4134.5454
$ws.Range("M141").Value = 1045.4546

$ws = $wb.Worksheets.Item("ARM")
# Row 32: set H=1528.4062, I=1416.4193, K=1416.4193, M=-1129.4193
$ws.Range("H32").Value = 1528.4062
$ws.Range("I32").Value = 1416.4193
$ws.Range("K32").Value = 1416.4193
$ws.Range("M32").Value = -1129.4193

# Row 61: set H=5999.25, J=6666, L=6666, N=-7090
$ws.Range("H61").Value = 5999.25
$ws.Range("J61").Value = 6666
$ws.Range("L61").Value = 6666
$ws.Range("N61").Value = -7090

# Row 132: set H=2671.4443, I=2671.4443, J=0, K=8014.3329, L=0; clear M
$ws.Range("H132").Value = 2671.4443
$ws.Range("I132").Value = 2671.4443
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 8014.3329
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()

# Row 136: set H=5999.25, J=6666, L=19998, N=-25098
$ws.Range("H136").Value = 5999.25
$ws.Range("J136").Value = 6666
$ws.Range("L136").Value = 19998
$ws.Range("N136").Value = -25098

$ws = $wb.Worksheets.Item("BSM")
# Row 86: set H=1476.2354, I=1584, J=1400.8, K=1584, L=1400.8, M=-461, N=-3646.8
$ws.Range("H86").Value = 1476.2354
$ws.Range("I86").Value = 1584
$ws.Range("J86").Value = 1400.8
$ws.Range("K86").Value = 1584
$ws.Range("L86").Value = 1400.8
$ws.Range("M86").Value = -461
$ws.Range("N86").Value = -3646.8

# Row 89: set H=1476.2354, I=1584, J=1400.8, K=7920, L=7004, M=-2304, N=-18236
$ws.Range("H89").Value = 1476.2354
$ws.Range("I89").Value = 1584
$ws.Range("J89").Value = 1400.8
$ws.Range("K89").Value = 7920
$ws.Range("L89").Value = 7004
$ws.Range("M89").Value = -2304
$ws.Range("N89").Value = -18236

# Row 94: set H=1571.8572, I=1290.6666, K=1290.6666, M=-839.6666
$ws.Range("H94").Value = 1571.8572
$ws.Range("I94").Value = 1290.6666
$ws.Range("K94").Value = 1290.6666
$ws.Range("M94").Value = -839.6666

# Row 107: set H=3131.3635, I=3244.5, K=3244.5, M=-1324.5
$ws.Range("H107").Value = 3131.3635
$ws.Range("I107").Value = 3244.5
$ws.Range("K107").Value = 3244.5
$ws.Range("M107").Value = -1324.5

# Row 134: set H=3361.1333, I=2529.7856, J=15000, K=7589.3568, L=45000, M=-5054.3568, N=-50070
$ws.Range("H134").Value = 3361.1333
$ws.Range("I134").Value = 2529.7856
$ws.Range("J134").Value = 15000
$ws.Range("K134").Value = 7589.3568
$ws.Range("L134").Value = 45000
$ws.Range("M134").Value = -5054.3568
$ws.Range("N134").Value = -50070

$ws = $wb.Worksheets.Item("CRP")
# Row 31: set H=7758.472, I=8667.214, J=7180.1816, K=8667.214, L=7180.1816, M=-8372.214, N=-7770.1816
$ws.Range("H31").Value = 7758.472
$ws.Range("I31").Value = 8667.214
$ws.Range("J31").Value = 7180.1816
$ws.Range("K31").Value = 8667.214
$ws.Range("L31").Value = 7180.1816
$ws.Range("M31").Value = -8372.214
$ws.Range("N31").Value = -7770.1816

# Row 34: set H=7758.472, I=8667.214, J=7180.1816, K=8667.214, L=7180.1816, M=-8465.214, N=-7584.1816
$ws.Range("H34").Value = 7758.472
$ws.Range("I34").Value = 8667.214
$ws.Range("J34").Value = 7180.1816
$ws.Range("K34").Value = 8667.214
$ws.Range("L34").Value = 7180.1816
$ws.Range("M34").Value = -8465.214
$ws.Range("N34").Value = -7584.1816

# Row 99: set H=8355.5, J=0, L=0
$ws.Range("H99").Value = 8355.5
$ws.Range("J99").Value = 0
$ws.Range("L99").Value = 0

# Row 105: set H=2591.7083, I=2899.5557, J=1668.1666, K=2899.5557, L=1668.1666, M=-1152.5557, N=-5162.1666
$ws.Range("H105").Value = 2591.7083
$ws.Range("I105").Value = 2899.5557
$ws.Range("J105").Value = 1668.1666
$ws.Range("K105").Value = 2899.5557
$ws.Range("L105").Value = 1668.1666
$ws.Range("M105").Value = -1152.5557
$ws.Range("N105").Value = -5162.1666

# Row 126: set H=8355.5, J=0, L=0
$ws.Range("H126").Value = 8355.5
$ws.Range("J126").Value = 0
$ws.Range("L126").Value = 0

# Row 132: set H=1725.5, I=1725.5, J=0, K=5176.5, L=0; clear M
$ws.Range("H132").Value = 1725.5
$ws.Range("I132").Value = 1725.5
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 5176.5
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
# Row 139: set H=94358.17999999999, I=127568.875, K=382706.625, M=-377566.625
$ws.Range("H139").Value = 94358.17999999999
$ws.Range("I139").Value = 127568.875
$ws.Range("K139").Value = 382706.625
$ws.Range("M139").Value = -377566.625

$ws = $wb.Worksheets.Item("GSM")
# Row 97: set H=1347.8334, I=1387, K=1387, M=-891
$ws.Range("H97").Value = 1347.8334
$ws.Range("I97").Value = 1387
$ws.Range("K97").Value = 1387
$ws.Range("M97").Value = -891

# Row 99: set H=4534.2, J=900, L=900, N=-5392
$ws.Range("H99").Value = 4534.2
$ws.Range("J99").Value = 900
$ws.Range("L99").Value = 900
$ws.Range("N99").Value = -5392

# Row 126: set H=4154.6924, I=3430.2856, J=4999.8335, K=10290.8568, L=14999.5005, M=-7820.856800000001, N=-19939.5005
$ws.Range("H126").Value = 4154.6924
$ws.Range("I126").Value = 3430.2856
$ws.Range("J126").Value = 4999.8335
$ws.Range("K126").Value = 10290.8568
$ws.Range("L126").Value = 14999.5005
$ws.Range("M126").Value = -7820.856800000001
$ws.Range("N126").Value = -19939.5005

# Row 132: set H=2082.0417, I=2082.0417, K=6246.125100000001, M=-3716.125100000001
$ws.Range("H132").Value = 2082.0417
$ws.Range("I132").Value = 2082.0417
$ws.Range("K132").Value = 6246.125100000001
$ws.Range("M132").Value = -3716.125100000001

$ws = $wb.Worksheets.Item("LTW")
# Row 46: set H=1739.091, I=1583, J=1797.625, K=1583, L=1797.625, M=-1395, N=-2173.625
$ws.Range("H46").Value = 1739.091
$ws.Range("I46").Value = 1583
$ws.Range("J46").Value = 1797.625
$ws.Range("K46").Value = 1583
$ws.Range("L46").Value = 1797.625
$ws.Range("M46").Value = -1395
$ws.Range("N46").Value = -2173.625

# Row 55: set H=291.06668, I=286.6875, J=296.07144, K=286.6875, L=296.07144, M=-113.6875, N=-642.0714399999999
$ws.Range("H55").Value = 291.06668
$ws.Range("I55").Value = 286.6875
$ws.Range("J55").Value = 296.07144
$ws.Range("K55").Value = 286.6875
$ws.Range("L55").Value = 296.07144
$ws.Range("M55").Value = -113.6875
$ws.Range("N55").Value = -642.0714399999999

# Row 68: set H=3548, I=2900, J=3936.8, K=2900, L=3936.8, M=-2151, N=-5434.8
$ws.Range("H68").Value = 3548
$ws.Range("I68").Value = 2900
$ws.Range("J68").Value = 3936.8
$ws.Range("K68").Value = 2900
$ws.Range("L68").Value = 3936.8
$ws.Range("M68").Value = -2151
$ws.Range("N68").Value = -5434.8

# Row 71: set H=3548, I=2900, J=3936.8, K=14500, L=19684, M=-10756, N=-27172
$ws.Range("H71").Value = 3548
$ws.Range("I71").Value = 2900
$ws.Range("J71").Value = 3936.8
$ws.Range("K71").Value = 14500
$ws.Range("L71").Value = 19684
$ws.Range("M71").Value = -10756
$ws.Range("N71").Value = -27172

# Row 82: set H=960.7143, I=864.8570999999999, J=1056.5714, K=864.8570999999999, L=1056.5714, M=-503.8570999999999, N=-1778.5714
$ws.Range("H82").Value = 960.7143
$ws.Range("I82").Value = 864.8570999999999
$ws.Range("J82").Value = 1056.5714
$ws.Range("K82").Value = 864.8570999999999
$ws.Range("L82").Value = 1056.5714
$ws.Range("M82").Value = -503.8570999999999
$ws.Range("N82").Value = -1778.5714

# Row 85: set H=960.7143, I=864.8570999999999, J=1056.5714, K=864.8570999999999, L=1056.5714, M=383.1429000000001, N=-3552.5714
$ws.Range("H85").Value = 960.7143
$ws.Range("I85").Value = 864.8570999999999
$ws.Range("J85").Value = 1056.5714
$ws.Range("K85").Value = 864.8570999999999
$ws.Range("L85").Value = 1056.5714
$ws.Range("M85").Value = 383.1429000000001
$ws.Range("N85").Value = -3552.5714

# Row 93: set H=3377.3, I=2971.75, J=4999.5, K=2971.75, L=4999.5, M=-1723.75, N=-7495.5
$ws.Range("H93").Value = 3377.3
$ws.Range("I93").Value = 2971.75
$ws.Range("J93").Value = 4999.5
$ws.Range("K93").Value = 2971.75
$ws.Range("L93").Value = 4999.5
$ws.Range("M93").Value = -1723.75
$ws.Range("N93").Value = -7495.5

# Row 99: set H=28612.666, I=28612.666, K=28612.666, M=-25617.666
$ws.Range("H99").Value = 28612.666
$ws.Range("I99").Value = 28612.666
$ws.Range("K99").Value = 28612.666
$ws.Range("M99").Value = -25617.666

$ws = $wb.Worksheets.Item("WVR")
# Row 81: set H=2809.9, I=2333.3333, J=3014.1428, K=4666.6666, L=6028.2856, M=-3605.6666, N=-8150.2856
$ws.Range("H81").Value = 2809.9
$ws.Range("I81").Value = 2333.3333
$ws.Range("J81").Value = 3014.1428
$ws.Range("K81").Value = 4666.6666
$ws.Range("L81").Value = 6028.2856
$ws.Range("M81").Value = -3605.6666
$ws.Range("N81").Value = -8150.2856

# Row 84: set H=2809.9, I=2333.3333, J=3014.1428, K=23333.333, L=30141.428, M=-18029.333, N=-40749.428
$ws.Range("H84").Value = 2809.9
$ws.Range("I84").Value = 2333.3333
$ws.Range("J84").Value = 3014.1428
$ws.Range("K84").Value = 23333.333
$ws.Range("L84").Value = 30141.428
$ws.Range("M84").Value = -18029.333
$ws.Range("N84").Value = -40749.428

# Row 107: set H=2005.6666, I=936.6429000000001, K=2809.9287, M=-889.9287000000004
$ws.Range("H107").Value = 2005.6666
$ws.Range("I107").Value = 936.6429000000001
$ws.Range("K107").Value = 2809.9287
$ws.Range("M107").Value = -889.9287000000004

# Row 126: set H=1610.4348, I=1522.1666, K=4566.4998, M=-2096.4998
$ws.Range("H126").Value = 1610.4348
$ws.Range("I126").Value = 1522.1666
$ws.Range("K126").Value = 4566.4998
$ws.Range("M126").Value = -2096.4998

# Row 132: set H=2786.8667, I=2469.0386, J=4852.75, K=7407.1158, L=14558.25, M=-4877.1158, N=-19618.25
$ws.Range("H132").Value = 2786.8667
$ws.Range("I132").Value = 2469.0386
$ws.Range("J132").Value = 4852.75
$ws.Range("K132").Value = 7407.1158
$ws.Range("L132").Value = 14558.25
$ws.Range("M132").Value = -4877.1158
$ws.Range("N132").Value = -19618.25
